# edit.ps1 — applies the LOM3107.docx content-shuffle edit described by the diff.
# Strategy: locate each affected paragraph by its (stable) 1-based Word Paragraphs
# index, then use Range.Find.Execute(old, ..., new, wdReplaceAll) scoped to that
# paragraph (or to a specific run-sized sub-range) so only the intended text is
# swapped while paragraph marks / pStyle / bold runs are left completely intact.

$d = $word.ActiveDocument

function Replace-InRange {
    param($range, $oldText, $newText)
    $ok = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find/Replace did not match expected text: $oldText"
    }
}

function Replace-ParagraphText {
    param($paraIndex, $oldText, $newText)
    $r = $d.Paragraphs.Item($paraIndex).Range
    Replace-InRange $r $oldText $newText
}

# --- Paragraph 6: "Objetivos" body text <-> "Programa resumido" summary text ---
$old6 = "Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões, deformações e deslocamentos em estruturas compostas por barras em regime elástico-linear sob carregamento axial, torção e flexão." + [char]11 + "Desenvolver aplicações práticas para dimensionamento de barras em condições de carregamentos mistos." + [char]11 + "Prover o conhecimento dos fenômenos de flambagem, com aplicações práticas para dimensionamento de colunas." + [char]11 + "Descrever a metodologia para análise dos estados planos de tensão e deformação, bem como a aplicação da lei de Hooke para casos multiaxiais." + [char]11 + "Apresentar conceitos básicos sobre energia de deformação."
$new6 = "Considerações fundamentais. Conceito de tensão. Conceito de deformação. Lei de Hooke. Carga Axial. Torção em barras de seção circular. Flexão em vigas isostáticas de seção simétrica. Cargas combinadas. Flambagem de colunas. Análise de Tensão e Deformação. Lei de Hooke Multiaxial. Energia de deformação."
Replace-ParagraphText 6 $old6 $new6

# --- Paragraph 8: "Docente(s)" list (4 runs, 1 break each) -> receives the former
# "Objetivos" sentences, the former "Programa" syllabus items, and the former
# "Metodo"/"Criterio" evaluation sentences -- one run at a time, in place, so the
# paragraph keeps exactly 4 runs (matching the target OOXML) instead of collapsing
# into one run. ---
$p8r0old = "471420 - Carlos Antonio Reis Pereira Baptista"
$p8r0new = "Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões, deformações e deslocamentos em estruturas compostas por barras em regime elástico-linear sob carregamento axial, torção e flexão." + [char]11 + "Desenvolver aplicações práticas para dimensionamento de barras em condições de carregamentos mistos." + [char]11 + "Prover o conhecimento dos fenômenos de flambagem, com aplicações práticas para dimensionamento de colunas." + [char]11 + "Descrever a metodologia para análise dos estados planos de tensão e deformação, bem como a aplicação da lei de Hooke para casos multiaxiais." + [char]11 + "Apresentar conceitos básicos sobre energia de deformação."
Replace-ParagraphText 8 $p8r0old $p8r0new
$p8r1old = "3480026 - João Paulo Pascon"
$p8r1new = "1. Considerações fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes." + [char]11 + "2. Conceito de tensão: Tensão Normal; Tensão Cisalhante; Tensões admissíveis." + [char]11 + "3. Conceito de deformação: Deformação Normal; Deformação por Cisalhamento." + [char]11 + "4. Lei de Hooke: Elasticidade linear e o Módulo de Young; Lei de Hooke para Cisalhamento." + [char]11 + "5. Carga Axial: Deslocamentos em sistemas isostáticos; Efeitos da Temperatura; Sistemas Hiperestáticos." + [char]11 + "6. Torção em barras de seção circular: Momento de inércia polar; Análise das tensões em eixos de seção maciça e seção vazada; Cálculo das rotações relativas entre seções adjacentes; Eixos estaticamente indeterminados; Torção e tração combinadas." + [char]11 + "7. Flexão em vigas isostáticas de seção simétrica: Forças concentradas e forças distribuídas; Diagramas de força cortante e momento fletor para uma viga carregada; Momento de inércia, eixos principais de inércia; Flexão em Vigas de Seção Simétrica; Determinação das Tensões Normais; Deflexões em vigas: equação diferencial da linha elástica; Tensões de cisalhamento em vigas. Tensões de cisalhamento em barras de paredes finas." + [char]11 + "8. Cargas combinadas: Modos Mistos de Carregamento. Projeto de barras submetidas a cargas axiais, transversais e torcionais." + [char]11 + "9. Flambagem de colunas: Raio de giração. Fórmula de Euler para colunas biarticuladas. Fatores de correção para outras condições de contorno. Projeto de colunas de aço e de outras ligas submetidas a um carregamento centrado." + [char]11 + "10. Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Transformação do Estado Plano de Deformação." + [char]11 + "11. Lei de Hooke Multiaxial: Elasticidade, Homogeneidade e Isotropia; Coeficiente de Poisson; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas." + [char]11 + "12. Energia de deformação: Densidade de energia de deformação. Energia de deformação elástica para tensões normais. Energia de deformação elástica para tensões de cisalhamento. Projeto para carregamento por impacto. Métodos de energia: teorema de Castigliano e suas aplicações."
Replace-ParagraphText 8 $p8r1old $p8r1new
$p8r2old = "5840793 - Sérgio Schneider"
$p8r2new = "Os alunos serão avaliados por meio de três conjuntos de notas: duas provas escritas (P1 e P2) envolvendo o conteúdo teórico ministrado em sala de aula; exercícios (EX) propostos regularmente para serem entregues e discutidos na aula subsequente; e seminários (SM) em grupo ao final da disciplina."
Replace-ParagraphText 8 $p8r2old $p8r2new
$p8r3old = "7797767 - Viktor Pastoukhov"
$p8r3new = "Nota Final (NF) = 70%((P1+P2)/2)+20%(EX)+10%(SM)."
Replace-ParagraphText 8 $p8r3old $p8r3new

# --- Paragraph 10: "Programa resumido" body text <-> "Norma de recuperacao" text ---
$old10 = "Considerações fundamentais. Conceito de tensão. Conceito de deformação. Lei de Hooke. Carga Axial. Torção em barras de seção circular. Flexão em vigas isostáticas de seção simétrica. Cargas combinadas. Flambagem de colunas. Análise de Tensão e Deformação. Lei de Hooke Multiaxial. Energia de deformação."
$new10 = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
Replace-ParagraphText 10 $old10 $new10

# --- Paragraph 12: "Programa" syllabus list <-> "Bibliografia" list ---
$old12 = "1. Considerações fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes." + [char]11 + "2. Conceito de tensão: Tensão Normal; Tensão Cisalhante; Tensões admissíveis." + [char]11 + "3. Conceito de deformação: Deformação Normal; Deformação por Cisalhamento." + [char]11 + "4. Lei de Hooke: Elasticidade linear e o Módulo de Young; Lei de Hooke para Cisalhamento." + [char]11 + "5. Carga Axial: Deslocamentos em sistemas isostáticos; Efeitos da Temperatura; Sistemas Hiperestáticos." + [char]11 + "6. Torção em barras de seção circular: Momento de inércia polar; Análise das tensões em eixos de seção maciça e seção vazada; Cálculo das rotações relativas entre seções adjacentes; Eixos estaticamente indeterminados; Torção e tração combinadas." + [char]11 + "7. Flexão em vigas isostáticas de seção simétrica: Forças concentradas e forças distribuídas; Diagramas de força cortante e momento fletor para uma viga carregada; Momento de inércia, eixos principais de inércia; Flexão em Vigas de Seção Simétrica; Determinação das Tensões Normais; Deflexões em vigas: equação diferencial da linha elástica; Tensões de cisalhamento em vigas. Tensões de cisalhamento em barras de paredes finas." + [char]11 + "8. Cargas combinadas: Modos Mistos de Carregamento. Projeto de barras submetidas a cargas axiais, transversais e torcionais." + [char]11 + "9. Flambagem de colunas: Raio de giração. Fórmula de Euler para colunas biarticuladas. Fatores de correção para outras condições de contorno. Projeto de colunas de aço e de outras ligas submetidas a um carregamento centrado." + [char]11 + "10. Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Transformação do Estado Plano de Deformação." + [char]11 + "11. Lei de Hooke Multiaxial: Elasticidade, Homogeneidade e Isotropia; Coeficiente de Poisson; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas." + [char]11 + "12. Energia de deformação: Densidade de energia de deformação. Energia de deformação elástica para tensões normais. Energia de deformação elástica para tensões de cisalhamento. Projeto para carregamento por impacto. Métodos de energia: teorema de Castigliano e suas aplicações."
$new12 = "1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p." + [char]11 + "2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p." + [char]11 + "3. R.R. CRAIG, Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p." + [char]11 + "4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p." + [char]11 + "5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p." + [char]11 + "6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p." + [char]11 + "7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p." + [char]11 + "8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais. Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p."
Replace-ParagraphText 12 $old12 $new12

# --- Paragraph 14: "Avaliacao" (Metodo:/Criterio:/Norma de recuperacao: bold labels)
# keep the bold label runs untouched; only the plain-text value runs after each
# label are swapped for the three docente names. ---
Replace-ParagraphText 14 "Os alunos serão avaliados por meio de três conjuntos de notas: duas provas escritas (P1 e P2) envolvendo o conteúdo teórico ministrado em sala de aula; exercícios (EX) propostos regularmente para serem entregues e discutidos na aula subsequente; e seminários (SM) em grupo ao final da disciplina." "471420 - Carlos Antonio Reis Pereira Baptista"
Replace-ParagraphText 14 "Nota Final (NF) = 70%((P1+P2)/2)+20%(EX)+10%(SM)." "3480026 - João Paulo Pascon"
Replace-ParagraphText 14 "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2." "5840793 - Sérgio Schneider"

# --- Paragraph 16: "Bibliografia" list <-> last docente name ---
$old16 = "1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p." + [char]11 + "2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p." + [char]11 + "3. R.R. CRAIG, Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p." + [char]11 + "4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p." + [char]11 + "5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p." + [char]11 + "6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p." + [char]11 + "7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p." + [char]11 + "8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais. Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p."
$new16 = "7797767 - Viktor Pastoukhov"
Replace-ParagraphText 16 $old16 $new16

Write-Output "LOM3107 content shuffle applied."
